# Refined metadata to be additional tab
#
# 1) Update the "time_taken" timestamps (column F) on the existing "data"
#    sheet for rows 2..17 (the panel query was re-run a bit later).
# 2) Add a new "metadata" worksheet (after "data") describing the panel
#    query itself (name/id/version/version-created/query-time/request url).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1) refresh data!F2:F17 timestamps -------------------------------------
$ws1.Range("F2").Value  = "2021-10-05 14:21:39.664560"
$ws1.Range("F3").Value  = "2021-10-05 14:21:39.664568"
$ws1.Range("F4").Value  = "2021-10-05 14:21:39.664572"
$ws1.Range("F5").Value  = "2021-10-05 14:21:39.664575"
$ws1.Range("F6").Value  = "2021-10-05 14:21:39.664578"
$ws1.Range("F7").Value  = "2021-10-05 14:21:39.664580"
$ws1.Range("F8").Value  = "2021-10-05 14:21:39.664583"
$ws1.Range("F9").Value  = "2021-10-05 14:21:39.664586"
$ws1.Range("F10").Value = "2021-10-05 14:21:39.664589"
$ws1.Range("F11").Value = "2021-10-05 14:21:39.664591"
$ws1.Range("F12").Value = "2021-10-05 14:21:39.664594"
$ws1.Range("F13").Value = "2021-10-05 14:21:39.664597"
$ws1.Range("F14").Value = "2021-10-05 14:21:39.664599"
$ws1.Range("F15").Value = "2021-10-05 14:21:39.664602"
$ws1.Range("F16").Value = "2021-10-05 14:21:39.664604"
$ws1.Range("F17").Value = "2021-10-05 14:21:39.664607"

# --- 2) add the "metadata" tab after "data" ---------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "metadata"

# Reuse the bold/bordered header style already used on the "data" sheet
# (copy format only, so no brand-new style entries are created).
$ws1.Range("B1").Copy()
$newSheet.Range("B1:G1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# header row
$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

# data row
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Multiple endocrine tumours"
$newSheet.Range("C2").Value = 36
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "1.13"
$newSheet.Range("E2").Value = "2021-08-02T08:17:08.350683Z"
$newSheet.Range("F2").Value = "2021-10-05 14:21:39.660893"
$newSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/36/?format=json"
